$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-09-05 Thursday", $false, $false, $false, $false, $false, $true, 1, $false, "2024-09-06 Friday", 2) | Out-Null

# Update the multiplication problems in the table (row, col -> new value)
$tbl = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="29×61=1769"},
    @{Row=1;  Col=2; Text="51×56=2856"},
    @{Row=1;  Col=3; Text="92×32=2944"},
    @{Row=1;  Col=4; Text="87×73=6351"},
    @{Row=1;  Col=5; Text="54×31=1674"},

    @{Row=5;  Col=1; Text="51×72=3672"},
    @{Row=5;  Col=2; Text="80×20=1600"},
    @{Row=5;  Col=3; Text="40×68=2720"},
    @{Row=5;  Col=4; Text="46×43=1978"},
    @{Row=5;  Col=5; Text="87×88=7656"},

    @{Row=10; Col=1; Text="95×11=1045"},
    @{Row=10; Col=2; Text="66×43=2838"},
    @{Row=10; Col=3; Text="31×17=527"},
    @{Row=10; Col=4; Text="80×81=6480"},
    @{Row=10; Col=5; Text="90×13=1170"},

    @{Row=15; Col=1; Text="60×33=1980"},
    @{Row=15; Col=2; Text="14×97=1358"},
    @{Row=15; Col=3; Text="98×96=9408"},
    @{Row=15; Col=4; Text="53×67=3551"},
    @{Row=15; Col=5; Text="64×47=3008"},

    @{Row=20; Col=1; Text="31×14=434"},
    @{Row=20; Col=2; Text="29×89=2581"},
    @{Row=20; Col=3; Text="29×88=2552"},
    @{Row=20; Col=4; Text="30×88=2640"},
    @{Row=20; Col=5; Text="58×95=5510"}
)

foreach ($u in $updates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $u.Text
}
